# Generate Report for Archive
# - Update localization status text "Ready for handoff" -> "In Translation"
#   on every sheet that surfaces it (Overview summary columns + per-language
#   status tables).
# - Shrink the "Status"/"zh-cn"/"de-de" status columns to match the new
#   (shorter) content width.

$wb = $excel.ActiveWorkbook

# --- Update the status text wherever it appears -----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Shrink the status columns to fit the new, shorter text -----------
$newColumnWidth = 13.4101845877511 - (5 / 6)

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth

$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth

$dede.Columns.Item(3).ColumnWidth = $newColumnWidth
